# Auto-generated edit script: applies the Faerie_Profits market-price refresh diff
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR (columns H-N on specific rows).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 588.5
$ws.Range("I2").Value = 282.85715
$ws.Range("J2").Value = 1301.6666
$ws.Range("K2").Value = 282.85715
$ws.Range("L2").Value = 1301.6666
$ws.Range("M2").Value = -169.85715
$ws.Range("N2").Value = -1527.6666
$ws.Range("H9").Value = 545
$ws.Range("I9").Value = 545
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 545
$ws.Range("L9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -376
$ws.Range("H11").Value = 96.25
$ws.Range("I11").Value = 96.25
$ws.Range("K11").Value = 96.25
$ws.Range("M11").Value = 43.75
$ws.Range("H19").Value = 899.5
$ws.Range("J19").Value = 1000
$ws.Range("L19").Value = 1000
$ws.Range("N19").Value = -1350
$ws.Range("H64").Value = 7022.2
$ws.Range("I64").Value = 4000
$ws.Range("J64").Value = 7777.75
$ws.Range("K64").Value = 4000
$ws.Range("L64").Value = 7777.75
$ws.Range("M64").Value = -3752
$ws.Range("N64").Value = -8273.75
$ws.Range("H67").Value = 7022.2
$ws.Range("I67").Value = 4000
$ws.Range("J67").Value = 7777.75
$ws.Range("K67").Value = 4000
$ws.Range("L67").Value = 7777.75
$ws.Range("M67").Value = -3142
$ws.Range("N67").Value = -9493.75
$ws.Range("H74").Value = 12460.467
$ws.Range("I74").Value = 13075.667
$ws.Range("K74").Value = 13075.667
$ws.Range("M74").Value = -12139.667
$ws.Range("H75").Value = 27642.5
$ws.Range("I75").Value = 20285
$ws.Range("J75").Value = 35000
$ws.Range("K75").Value = 20285
$ws.Range("L75").Value = 35000
$ws.Range("M75").Value = -19349
$ws.Range("N75").Value = -36872
$ws.Range("H77").Value = 12460.467
$ws.Range("I77").Value = 13075.667
$ws.Range("K77").Value = 65378.335
$ws.Range("M77").Value = -60698.335
$ws.Range("H78").Value = 27642.5
$ws.Range("I78").Value = 20285
$ws.Range("J78").Value = 35000
$ws.Range("K78").Value = 60855
$ws.Range("L78").Value = 105000
$ws.Range("M78").Value = -56175
$ws.Range("N78").Value = -114360
$ws.Range("H88").Value = 1879.9375
$ws.Range("I88").Value = 1279.6
$ws.Range("K88").Value = 1279.6
$ws.Range("M88").Value = -873.5999999999999
$ws.Range("H91").Value = 1879.9375
$ws.Range("I91").Value = 1279.6
$ws.Range("K91").Value = 1279.6
$ws.Range("M91").Value = 124.4000000000001
$ws.Range("H96").Value = 1714
$ws.Range("I96").Value = 1801.7142
$ws.Range("J96").Value = 1100
$ws.Range("K96").Value = 5405.142599999999
$ws.Range("L96").Value = 3300
$ws.Range("M96").Value = -4032.142599999999
$ws.Range("N96").Value = -6046
$ws.Range("H100").Value = 9756.5
$ws.Range("I100").Value = 1536.6666
$ws.Range("J100").Value = 13279.286
$ws.Range("K100").Value = 1536.6666
$ws.Range("L100").Value = 13279.286
$ws.Range("M100").Value = -995.6666
$ws.Range("N100").Value = -14361.286
$ws.Range("H115").Value = 934.3333
$ws.Range("I115").Value = 934.3333
$ws.Range("K115").Value = 2802.9999
$ws.Range("M115").Value = -1235.9999
$ws.Range("H125").Value = 4655.4443
$ws.Range("J125").Value = 5112.375
$ws.Range("L125").Value = 46011.375
$ws.Range("N125").Value = -50931.375
$ws.Range("H129").Value = 52634730
$ws.Range("I129").Value = 90909540
$ws.Range("K129").Value = 272728620
$ws.Range("M129").Value = -272723620
$ws.Range("H130").Value = 81900
$ws.Range("I130").Value = 60350
$ws.Range("J130").Value = 125000
$ws.Range("K130").Value = 60350
$ws.Range("L130").Value = 125000
$ws.Range("M130").Value = -55330
$ws.Range("N130").Value = -135040
$ws.Range("H132").Value = 23817754
$ws.Range("I132").Value = 30305856
$ws.Range("K132").Value = 90917568
$ws.Range("M132").Value = -90915038
$ws.Range("H135").Value = 1667.16
$ws.Range("I135").Value = 482.2381
$ws.Range("K135").Value = 4340.1429
$ws.Range("M135").Value = -1805.1429
$ws.Range("H137").Value = 1793.8529
$ws.Range("I137").Value = 1857
$ws.Range("J137").Value = 1678.0834
$ws.Range("K137").Value = 5571
$ws.Range("L137").Value = 5034.2502
$ws.Range("M137").Value = -3021
$ws.Range("N137").Value = -10134.2502
$ws.Range("H138").Value = 356488.1
$ws.Range("I138").Value = 92200.73
$ws.Range("J138").Value = 501846.16
$ws.Range("K138").Value = 276602.19
$ws.Range("L138").Value = 1505538.48
$ws.Range("M138").Value = -271462.19
$ws.Range("N138").Value = -1515818.48

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9767.784
$ws.Range("I32").Value = 10345.879
$ws.Range("J32").Value = 4998.5
$ws.Range("K32").Value = 10345.879
$ws.Range("L32").Value = 4998.5
$ws.Range("M32").Value = -10058.879
$ws.Range("N32").Value = -5572.5
$ws.Range("H61").Value = 4856.3433
$ws.Range("I61").Value = 3167.9778
$ws.Range("K61").Value = 3167.9778
$ws.Range("M61").Value = -2955.9778
$ws.Range("H74").Value = 3099.158
$ws.Range("I74").Value = 2489.9062
$ws.Range("K74").Value = 2489.9062
$ws.Range("M74").Value = -1615.9062
$ws.Range("H76").Value = 43507.332
$ws.Range("I76").Value = 35261
$ws.Range("J76").Value = 60000
$ws.Range("K76").Value = 35261
$ws.Range("L76").Value = 60000
$ws.Range("M76").Value = -34923
$ws.Range("N76").Value = -60676
$ws.Range("H77").Value = 3099.158
$ws.Range("I77").Value = 2489.9062
$ws.Range("K77").Value = 12449.531
$ws.Range("M77").Value = -8081.530999999999
$ws.Range("H79").Value = 43507.332
$ws.Range("I79").Value = 35261
$ws.Range("J79").Value = 60000
$ws.Range("K79").Value = 35261
$ws.Range("L79").Value = 60000
$ws.Range("M79").Value = -34091
$ws.Range("N79").Value = -62340
$ws.Range("H97").Value = 1000
$ws.Range("I97").Value = 1000
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1000
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -504
$ws.Range("H102").Value = 22224920
$ws.Range("I102").Value = 2599
$ws.Range("K102").Value = 2599
$ws.Range("M102").Value = -977
$ws.Range("H110").Value = 2344
$ws.Range("I110").Value = 2388.4
$ws.Range("J110").Value = 1900
$ws.Range("K110").Value = 2388.4
$ws.Range("L110").Value = 1900
$ws.Range("M110").Value = -343.4000000000001
$ws.Range("N110").Value = -5990
$ws.Range("H132").Value = 2456.0364
$ws.Range("I132").Value = 2333.82
$ws.Range("J132").Value = 3678.2
$ws.Range("K132").Value = 7001.460000000001
$ws.Range("L132").Value = 11034.6
$ws.Range("M132").Value = -4471.460000000001
$ws.Range("N132").Value = -16094.6
$ws.Range("H136").Value = 4856.3433
$ws.Range("I136").Value = 3167.9778
$ws.Range("K136").Value = 9503.9334
$ws.Range("M136").Value = -6953.9334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2790.7273
$ws.Range("I20").Value = 2601.9167
$ws.Range("K20").Value = 2601.9167
$ws.Range("M20").Value = -2354.9167
$ws.Range("H57").Value = 46996.668
$ws.Range("J57").Value = 47995
$ws.Range("L57").Value = 47995
$ws.Range("N57").Value = -49435
$ws.Range("H94").Value = 1455.3334
$ws.Range("I94").Value = 1599
$ws.Range("K94").Value = 1599
$ws.Range("M94").Value = -1148
$ws.Range("H99").Value = 3460.5652
$ws.Range("I99").Value = 2685.2856
$ws.Range("K99").Value = 2685.2856
$ws.Range("M99").Value = -1187.2856
$ws.Range("H105").Value = 6239.9
$ws.Range("I105").Value = 3733.1667
$ws.Range("J105").Value = 10000
$ws.Range("K105").Value = 3733.1667
$ws.Range("L105").Value = 10000
$ws.Range("M105").Value = -1986.1667
$ws.Range("N105").Value = -13494
$ws.Range("H107").Value = 2458.4
$ws.Range("I107").Value = 2447.5
$ws.Range("J107").Value = 2465.6667
$ws.Range("K107").Value = 2447.5
$ws.Range("L107").Value = 2465.6667
$ws.Range("M107").Value = -527.5
$ws.Range("N107").Value = -6305.6667
$ws.Range("H131").Value = 110000
$ws.Range("J131").Value = 110000
$ws.Range("L131").Value = 110000
$ws.Range("N131").Value = -120080
$ws.Range("H134").Value = 8611.474
$ws.Range("I134").Value = 4319.5713
$ws.Range("K134").Value = 12958.7139
$ws.Range("M134").Value = -10423.7139
$ws.Range("H136").Value = 46996.668
$ws.Range("J136").Value = 47995
$ws.Range("L136").Value = 47995
$ws.Range("N136").Value = -58195

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2145.6428
$ws.Range("I58").Value = 2320.125
$ws.Range("J58").Value = 1913
$ws.Range("K58").Value = 2320.125
$ws.Range("L58").Value = 1913
$ws.Range("M58").Value = -2117.125
$ws.Range("N58").Value = -2319
$ws.Range("H107").Value = 1416.2727
$ws.Range("J107").Value = 3137.6667
$ws.Range("L107").Value = 3137.6667
$ws.Range("N107").Value = -6977.6667
$ws.Range("H122").Value = 3772.524
$ws.Range("I122").Value = 4158
$ws.Range("K122").Value = 12474
$ws.Range("M122").Value = -10024
$ws.Range("H125").Value = 90001
$ws.Range("J125").Value = 90001
$ws.Range("L125").Value = 90001
$ws.Range("N125").Value = -94921
$ws.Range("H132").Value = 1484335.6
$ws.Range("I132").Value = 1669211.2
$ws.Range("K132").Value = 5007633.6
$ws.Range("M132").Value = -5005103.6
$ws.Range("H134").Value = 2792.9553
$ws.Range("I134").Value = 1247.551
$ws.Range("K134").Value = 3742.653
$ws.Range("M134").Value = -1207.653
$ws.Range("H136").Value = 2145.6428
$ws.Range("I136").Value = 2320.125
$ws.Range("J136").Value = 1913
$ws.Range("K136").Value = 6960.375
$ws.Range("L136").Value = 5739
$ws.Range("M136").Value = -4410.375
$ws.Range("N136").Value = -10839
$ws.Range("H141").Value = 368854.12
$ws.Range("J141").Value = 368854.12
$ws.Range("L141").Value = 368854.12
$ws.Range("N141").Value = -379214.12

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 973.2
$ws.Range("I5").Value = 506
$ws.Range("J5").Value = 1674
$ws.Range("K5").Value = 1518
$ws.Range("L5").Value = 5022
$ws.Range("M5").Value = -1406
$ws.Range("N5").Value = -5246
$ws.Range("H63").Value = 5915.5713
$ws.Range("I63").Value = 4705.5
$ws.Range("K63").Value = 14116.5
$ws.Range("M63").Value = -13367.5
$ws.Range("H66").Value = 5915.5713
$ws.Range("I66").Value = 4705.5
$ws.Range("K66").Value = 42349.5
$ws.Range("M66").Value = -38605.5
$ws.Range("H103").Value = 567.875
$ws.Range("J103").Value = 621.6667
$ws.Range("L103").Value = 1865.0001
$ws.Range("N103").Value = -3623.0001
$ws.Range("H107").Value = 4375.6924
$ws.Range("I107").Value = 386.25
$ws.Range("K107").Value = 1158.75
$ws.Range("M107").Value = 761.25
$ws.Range("H133").Value = 5442.25
$ws.Range("J133").Value = 6589.6665
$ws.Range("L133").Value = 19768.9995
$ws.Range("N133").Value = -29888.9995
$ws.Range("H135").Value = 973.2
$ws.Range("I135").Value = 506
$ws.Range("J135").Value = 1674
$ws.Range("K135").Value = 4554
$ws.Range("L135").Value = 15066
$ws.Range("M135").Value = -2019
$ws.Range("N135").Value = -20136
$ws.Range("H136").Value = 3159.4666
$ws.Range("I136").Value = 2324.3333
$ws.Range("K136").Value = 6972.999899999999
$ws.Range("M136").Value = -1872.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2845.1765
$ws.Range("I80").Value = 2696.5833
$ws.Range("J80").Value = 3201.8
$ws.Range("K80").Value = 2696.5833
$ws.Range("L80").Value = 3201.8
$ws.Range("M80").Value = -1698.5833
$ws.Range("N80").Value = -5197.8
$ws.Range("H83").Value = 2845.1765
$ws.Range("I83").Value = 2696.5833
$ws.Range("J83").Value = 3201.8
$ws.Range("K83").Value = 13482.9165
$ws.Range("L83").Value = 16009
$ws.Range("M83").Value = -8490.916499999999
$ws.Range("N83").Value = -25993
$ws.Range("H97").Value = 2809
$ws.Range("I97").Value = 1720
$ws.Range("J97").Value = 4170.25
$ws.Range("K97").Value = 1720
$ws.Range("L97").Value = 4170.25
$ws.Range("M97").Value = -1224
$ws.Range("N97").Value = -5162.25
$ws.Range("H98").Value = 28500
$ws.Range("J98").Value = 28500
$ws.Range("L98").Value = 28500
$ws.Range("N98").Value = -34490
$ws.Range("H105").Value = 250037250
$ws.Range("J105").Value = 250037250
$ws.Range("L105").Value = 250037250
$ws.Range("N105").Value = -250044238
$ws.Range("H106").Value = 23280
$ws.Range("J106").Value = 23280
$ws.Range("L106").Value = 23280
$ws.Range("N106").Value = -25804
$ws.Range("H107").Value = 843.4666999999999
$ws.Range("J107").Value = 1294
$ws.Range("L107").Value = 1294
$ws.Range("N107").Value = -5134
$ws.Range("H113").Value = 2242.5715
$ws.Range("I113").Value = 2249.6667
$ws.Range("J113").Value = 2200
$ws.Range("K113").Value = 2249.6667
$ws.Range("L113").Value = 2200
$ws.Range("M113").Value = -79.66670000000022
$ws.Range("N113").Value = -6540
$ws.Range("H122").Value = 1706.4615
$ws.Range("I122").Value = 1327.1904
$ws.Range("J122").Value = 3299.4
$ws.Range("K122").Value = 3981.5712
$ws.Range("L122").Value = 9898.200000000001
$ws.Range("M122").Value = -1531.5712
$ws.Range("N122").Value = -14798.2
$ws.Range("H126").Value = 2698.5
$ws.Range("I126").Value = 2266.3333
$ws.Range("K126").Value = 6798.999899999999
$ws.Range("M126").Value = -4328.999899999999
$ws.Range("H132").Value = 6291718
$ws.Range("I132").Value = 7409502
$ws.Range("K132").Value = 22228506
$ws.Range("M132").Value = -22225976

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3090.7273
$ws.Range("I22").Value = 745
$ws.Range("J22").Value = 3325.3
$ws.Range("K22").Value = 745
$ws.Range("L22").Value = 3325.3
$ws.Range("M22").Value = -450
$ws.Range("N22").Value = -3915.3
$ws.Range("H27").Value = 3090.7273
$ws.Range("I27").Value = 745
$ws.Range("J27").Value = 3325.3
$ws.Range("K27").Value = 745
$ws.Range("L27").Value = 3325.3
$ws.Range("M27").Value = -638
$ws.Range("N27").Value = -3539.3
$ws.Range("H38").Value = 5303.0303
$ws.Range("I38").Value = 5000
$ws.Range("J38").Value = 10000
$ws.Range("K38").Value = 5000
$ws.Range("L38").Value = 10000
$ws.Range("M38").Value = -4590
$ws.Range("N38").Value = -10820
$ws.Range("H40").Value = 4501.1953
$ws.Range("I40").Value = 4554.5625
$ws.Range("J40").Value = 4311.4443
$ws.Range("K40").Value = 4554.5625
$ws.Range("L40").Value = 4311.4443
$ws.Range("M40").Value = -4418.5625
$ws.Range("N40").Value = -4583.4443
$ws.Range("H55").Value = 268
$ws.Range("J55").Value = 253
$ws.Range("L55").Value = 253
$ws.Range("N55").Value = -599
$ws.Range("H68").Value = 4222.222
$ws.Range("I68").Value = 3916.5
$ws.Range("J68").Value = 4833.6665
$ws.Range("K68").Value = 3916.5
$ws.Range("L68").Value = 4833.6665
$ws.Range("M68").Value = -3167.5
$ws.Range("N68").Value = -6331.6665
$ws.Range("H71").Value = 4222.222
$ws.Range("I71").Value = 3916.5
$ws.Range("J71").Value = 4833.6665
$ws.Range("K71").Value = 19582.5
$ws.Range("L71").Value = 24168.3325
$ws.Range("M71").Value = -15838.5
$ws.Range("N71").Value = -31656.3325
$ws.Range("H82").Value = 8538.727999999999
$ws.Range("I82").Value = 9922.615
$ws.Range("K82").Value = 9922.615
$ws.Range("M82").Value = -9561.615
$ws.Range("H85").Value = 8538.727999999999
$ws.Range("I85").Value = 9922.615
$ws.Range("K85").Value = 9922.615
$ws.Range("M85").Value = -8674.615
$ws.Range("H93").Value = 3348
$ws.Range("I93").Value = 2464
$ws.Range("K93").Value = 2464
$ws.Range("M93").Value = -1216
$ws.Range("H95").Value = 48926.2
$ws.Range("J95").Value = 48926.2
$ws.Range("L95").Value = 48926.2
$ws.Range("N95").Value = -54418.2
$ws.Range("H98").Value = 72399.5
$ws.Range("J98").Value = 72399.5
$ws.Range("L98").Value = 72399.5
$ws.Range("N98").Value = -78389.5
$ws.Range("H100").Value = 3851.0527
$ws.Range("I100").Value = 3321.5386
$ws.Range("J100").Value = 4998.3335
$ws.Range("K100").Value = 3321.5386
$ws.Range("L100").Value = 4998.3335
$ws.Range("M100").Value = -2780.5386
$ws.Range("N100").Value = -6080.3335
$ws.Range("H132").Value = 3113.848
$ws.Range("I132").Value = 3111.3057
$ws.Range("K132").Value = 9333.917099999999
$ws.Range("M132").Value = -6803.917099999999
$ws.Range("H136").Value = 4383.879
$ws.Range("I136").Value = 4171.923
$ws.Range("J136").Value = 5171.143
$ws.Range("K136").Value = 12515.769
$ws.Range("L136").Value = 15513.429
$ws.Range("M136").Value = -9965.769
$ws.Range("N136").Value = -20613.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 159530.2
$ws.Range("I62").Value = 198412.75
$ws.Range("K62").Value = 198412.75
$ws.Range("M62").Value = -197788.75
$ws.Range("H65").Value = 159530.2
$ws.Range("I65").Value = 198412.75
$ws.Range("K65").Value = 992063.75
$ws.Range("M65").Value = -988943.75
$ws.Range("H81").Value = 55733.75
$ws.Range("I81").Value = 86443.164
$ws.Range("K81").Value = 172886.328
$ws.Range("M81").Value = -171825.328
$ws.Range("H84").Value = 55733.75
$ws.Range("I84").Value = 86443.164
$ws.Range("K84").Value = 864431.64
$ws.Range("M84").Value = -859127.64
$ws.Range("H96").Value = 3517.875
$ws.Range("I96").Value = 2264.3333
$ws.Range("J96").Value = 4270
$ws.Range("K96").Value = 2264.3333
$ws.Range("L96").Value = 4270
$ws.Range("M96").Value = -891.3332999999998
$ws.Range("N96").Value = -7016
$ws.Range("H100").Value = 3385.7778
$ws.Range("I100").Value = 2579.3333
$ws.Range("K100").Value = 5158.6666
$ws.Range("M100").Value = -4617.6666
$ws.Range("H101").Value = 68533.664
$ws.Range("J101").Value = 68533.664
$ws.Range("L101").Value = 68533.664
$ws.Range("N101").Value = -75023.664
$ws.Range("H107").Value = 1398.8572
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1398.8572
$ws.Range("K107").Value = 0
$ws.Range("L107").ClearContents()
$ws.Range("M107").Value = 4196.571599999999
$ws.Range("N107").Value = -8036.571599999999
$ws.Range("H122").Value = 2224.35
$ws.Range("I122").Value = 1945.7646
$ws.Range("K122").Value = 5837.293799999999
$ws.Range("M122").Value = -3387.293799999999
$ws.Range("H132").Value = 5062.364
$ws.Range("I132").Value = 7247.5
$ws.Range("J132").Value = 2440.2
$ws.Range("K132").Value = 21742.5
$ws.Range("L132").Value = 7320.599999999999
$ws.Range("M132").Value = -19212.5
$ws.Range("N132").Value = -12380.6
$ws.Range("H136").Value = 21129.8
$ws.Range("I136").Value = 26162.25
$ws.Range("K136").Value = 78486.75
$ws.Range("M136").Value = -75936.75
$ws.Range("H141").Value = 87530
$ws.Range("J141").Value = 82250
$ws.Range("L141").Value = 82250
$ws.Range("N141").Value = -92610

